$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.024194244416059
$ws.Range("D2").Value = 1.029408791303877
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.035628406177603
$ws.Range("I2").Value = 1.034190942919602
$ws.Range("J2").Value = 1.029370224458952
$ws.Range("K2").Value = 1.032222787573555
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.038424465723949
$ws.Range("N2").Value = 1.013858947106946
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.02493184508438
$ws.Range("D3").Value = 1.029947926262335
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.036649747987729
$ws.Range("I3").Value = 1.03436346511185
$ws.Range("J3").Value = 1.029747723985935
$ws.Range("K3").Value = 1.032570786295728
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.039254678676631
$ws.Range("N3").Value = 1.013983417266436
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.025409447797493
$ws.Range("D4").Value = 1.030296983803465
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.03731131940605
$ws.Range("I4").Value = 1.034473883859476
$ws.Range("J4").Value = 1.029991630085582
$ws.Range("K4").Value = 1.032795456917344
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.039791972428522
$ws.Range("N4").Value = 1.014063833146252
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.025610308282833
$ws.Range("D5").Value = 1.030443774057773
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.037589609196631
$ws.Range("I5").Value = 1.034520012365098
$ws.Range("J5").Value = 1.030094080613452
$ws.Range("K5").Value = 1.032889785739416
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.040017871217815
$ws.Range("N5").Value = 1.014097609768685
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.025644038067763
$ws.Range("D6").Value = 1.03046842346261
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.037636344927298
$ws.Range("I6").Value = 1.034527740426648
$ws.Range("J6").Value = 1.030111277339418
$ws.Range("K6").Value = 1.032905616735091
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.040055801786781
$ws.Range("N6").Value = 1.014103279229561
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.025412131407414
$ws.Range("D7").Value = 1.030298945041597
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.037315037280373
$ws.Range("I7").Value = 1.034474501377297
$ws.Range("J7").Value = 1.029992999380114
$ws.Range("K7").Value = 1.032796717827262
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.039794990820115
$ws.Range("N7").Value = 1.014064284590324
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.024443451553522
$ws.Range("D8").Value = 1.029590951814497
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.035973428796413
$ws.Range("I8").Value = 1.034249498653834
$ws.Range("J8").Value = 1.029497876352954
$ws.Range("K8").Value = 1.032340499837887
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.038705020650847
$ws.Range("N8").Value = 1.013901037961541
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.022739083122468
$ws.Range("D9").Value = 1.028344992864498
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.033614719364497
$ws.Range("I9").Value = 1.033843745342608
$ws.Range("J9").Value = 1.028622687069982
$ws.Range("K9").Value = 1.031532746153553
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.036785097918241
$ws.Range("N9").Value = 1.013612438527129
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.021604664830613
$ws.Range("D10").Value = 1.027515540668235
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.032045931871096
$ws.Range("I10").Value = 1.033567054967109
$ws.Range("J10").Value = 1.028037464023716
$ws.Range("K10").Value = 1.030991729650614
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.035505717607282
$ws.Range("N10").Value = 1.013419430487111
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.021113903593817
$ws.Range("D11").Value = 1.027156681560487
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.031367518247311
$ws.Range("I11").Value = 1.0334457865026
$ws.Range("J11").Value = 1.027783651292872
$ws.Range("K11").Value = 1.030756881375027
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.034951880569993
$ws.Range("N11").Value = 1.013335716240362
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.020931682214228
$ws.Range("D12").Value = 1.027023431976821
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.03111565881194
$ws.Range("I12").Value = 1.033400523415708
$ws.Range("J12").Value = 1.027689313836806
$ws.Range("K12").Value = 1.030669561573957
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.034746183413685
$ws.Range("N12").Value = 1.013304600307503
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.020970766160503
$ws.Range("D13").Value = 1.027052012325489
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.031169677452439
$ws.Range("I13").Value = 1.033410242381052
$ws.Range("J13").Value = 1.027709552243946
$ws.Range("K13").Value = 1.030688295872732
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.034790305112809
$ws.Range("N13").Value = 1.013311275713014
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.021098839696456
$ws.Range("D14").Value = 1.027145666146408
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.031346696736091
$ws.Range("I14").Value = 1.033442049496486
$ws.Range("J14").Value = 1.027775854549225
$ws.Range("K14").Value = 1.030749665255062
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.034934877112009
$ws.Range("N14").Value = 1.013333144608463
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.021177759248215
$ws.Range("D15").Value = 1.027203375591662
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.031455781894882
$ws.Range("I15").Value = 1.03346161796782
$ws.Range("J15").Value = 1.027816697655412
$ws.Range("K15").Value = 1.030787465502762
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.035023955716563
$ws.Range("N15").Value = 1.013346616020371
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.02163724464893
$ws.Range("D16").Value = 1.027539363415455
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.032090974621496
$ws.Range("I16").Value = 1.033575072463234
$ws.Range("J16").Value = 1.028054300250181
$ws.Range("K16").Value = 1.031007303537601
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.035542477054318
$ws.Range("N16").Value = 1.013424983397823
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.021925589198184
$ws.Range("D17").Value = 1.027750201012609
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.032489651088956
$ws.Range("I17").Value = 1.033645849177535
$ws.Range("J17").Value = 1.028203233803158
$ws.Range("K17").Value = 1.031145046423057
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.035867770860453
$ws.Range("N17").Value = 1.013474103833236
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.022093818940265
$ws.Range("D18").Value = 1.027873207766119
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.03272227732436
$ws.Range("I18").Value = 1.033686991148477
$ws.Range("J18").Value = 1.028290064854703
$ws.Range("K18").Value = 1.03122533304064
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.03605752296641
$ws.Range("N18").Value = 1.013502741371771
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.022151188218758
$ws.Range("D19").Value = 1.02791515471438
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.032801611241876
$ws.Range("I19").Value = 1.03370099557809
$ws.Range("J19").Value = 1.028319665288499
$ws.Range("K19").Value = 1.03125269908858
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.036122225827258
$ws.Range("N19").Value = 1.013512503717519
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.021894648071743
$ws.Range("D20").Value = 1.027727577147589
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.032446868077242
$ws.Range("I20").Value = 1.033638270073497
$ws.Range("J20").Value = 1.028187258706826
$ws.Range("K20").Value = 1.031130273748569
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.035832868485347
$ws.Range("N20").Value = 1.013468835077375
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.021061123308644
$ws.Range("D21").Value = 1.027118086128375
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.031294565272342
$ws.Range("I21").Value = 1.033432689121812
$ws.Range("J21").Value = 1.027756331823409
$ws.Range("K21").Value = 1.03073159588571
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.034892303629608
$ws.Range("N21").Value = 1.013326705333196
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.020537454380015
$ws.Range("D22").Value = 1.026735145783327
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.030570840238699
$ws.Range("I22").Value = 1.033302167942113
$ws.Range("J22").Value = 1.027485043764079
$ws.Range("K22").Value = 1.03048043025685
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.034301064346837
$ws.Range("N22").Value = 1.013237222944088
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.020815022599233
$ws.Range("D23").Value = 1.026938123459251
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.030954426879778
$ws.Range("I23").Value = 1.033371479285755
$ws.Range("J23").Value = 1.027628891248667
$ws.Range("K23").Value = 1.030613625024833
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.034614478681251
$ws.Range("N23").Value = 1.01328467048081
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.021908628901521
$ws.Range("D24").Value = 1.027737799809746
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.032466199616112
$ws.Range("I24").Value = 1.033641695179779
$ws.Range("J24").Value = 1.028194477286996
$ws.Range("K24").Value = 1.031136949058929
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.03584863932395
$ws.Range("N24").Value = 1.013471215843364
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.023179388447946
$ws.Range("D25").Value = 1.028666901249198
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.034223857902699
$ws.Range("I25").Value = 1.03394973624973
$ws.Range("J25").Value = 1.028849260044448
$ws.Range("K25").Value = 1.03174201787773
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.013687157262542

Write-Host "Applied vm_pu updates for 380 kV case"
